$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '61.582.69'
$ws.Cells.Item(2, 5).Value = '  -1.98%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.029.09'
$ws.Cells.Item(3, 5).Value = '  -1.31%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.01%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '535.79'
$ws.Cells.Item(5, 5).Value = '  -0.17%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '136.10'
$ws.Cells.Item(6, 5).Value = '  +1.67%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.03%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '3.019.91'
$ws.Cells.Item(8, 5).Value = '  -1.35%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.499'
$ws.Cells.Item(9, 5).Value = '  +1.21%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.150'
$ws.Cells.Item(10, 5).Value = '  -2.31%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '6.16'
$ws.Cells.Item(11, 5).Value = '  +0.47%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.452'
$ws.Cells.Item(12, 5).Value = '  -0.28%  '

# Row 13
$ws.Cells.Item(13, 2).Value = 'ShibaInu'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.0000222'
$ws.Cells.Item(13, 5).Value = '  -0.43%  '

# Row 14
$ws.Cells.Item(14, 2).Value = 'Avalanche'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '34.68'
$ws.Cells.Item(14, 5).Value = '  +1.42%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '3.518.78'
$ws.Cells.Item(15, 5).Value = '  -1.07%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  +0.39%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '61.595.81'
$ws.Cells.Item(17, 5).Value = '  -2.00%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '3.031.86'
$ws.Cells.Item(18, 5).Value = '  -0.83%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '6.69'
$ws.Cells.Item(19, 5).Value = '  +0.92%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '471.23'
$ws.Cells.Item(20, 5).Value = '  -2.09%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '13.39'
$ws.Cells.Item(21, 5).Value = '  +0.65%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.686'
$ws.Cells.Item(22, 5).Value = '  -0.98%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '7.04'
$ws.Cells.Item(23, 5).Value = '  -1.17%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '80.00'
$ws.Cells.Item(24, 5).Value = '  +1.23%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '12.24'
$ws.Cells.Item(25, 5).Value = '  +1.16%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.08%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '2.71'
$ws.Cells.Item(27, 5).Value = '  +0.60%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '7.94'
$ws.Cells.Item(28, 5).Value = '  -1.78%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.00'
$ws.Cells.Item(29, 5).Value = '  +0.42%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.93'
$ws.Cells.Item(30, 5).Value = '  +3.90%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '25.89'
$ws.Cells.Item(31, 5).Value = '  -0.38%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +2.66%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '5.60'
$ws.Cells.Item(33, 5).Value = '  +4.71%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '2.31'
$ws.Cells.Item(34, 5).Value = '  -1.58%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '55.61'
$ws.Cells.Item(35, 5).Value = '  -1.82%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '5.97'
$ws.Cells.Item(36, 5).Value = '  +0.02%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '470.37'
$ws.Cells.Item(37, 5).Value = '  -1.39%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '3.221.45'
$ws.Cells.Item(38, 5).Value = '  +4.22%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'Hedera'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.0799'
$ws.Cells.Item(39, 5).Value = '  +0.49%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'VeChain'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.0391'
$ws.Cells.Item(40, 5).Value = '  -0.75%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +3.38%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '8.25'
$ws.Cells.Item(42, 5).Value = '  +1.88%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '2.53'
$ws.Cells.Item(43, 5).Value = '  -4.28%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '27.96'
$ws.Cells.Item(44, 5).Value = '  +14.48%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.250'
$ws.Cells.Item(45, 5).Value = '  -0.42%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  +1.71%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +1.26%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '119.38'
$ws.Cells.Item(49, 5).Value = '  -1.47%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '0.0₃0504'
$ws.Cells.Item(50, 5).Value = '  -6.61%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +8.43%  '
